$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The raw-views export now writes the film_id column as a genuine number
# instead of an inline string, and the filtered rows for the old
# 10x0x0x ids (which had no data beyond column A/D/E placeholders) were
# dropped from the refreshed source, shifting the later "placeholder"
# rows up.

# 1) Drop the 18 now-stale placeholder rows (old rows 151-168).
$ws.Rows("151:168").Delete()

# 2) Re-write column A (film_id) for every remaining data row as a true
#    number rather than a text string.
$filmIds = @(1070101, 1070102, 1070103, 1070105, 1070106, 1070107, 1070109, 1070110, 1070111, 1070112, 1070113, 1070114, 1070115, 1070117, 1080101, 1080102, 1080103, 1080104, 1080105, 1080106, 1080107, 1080108, 1080109, 1080110, 1080111, 1080112, 1080113, 1080114, 1080201, 1080202, 1080203, 1080301, 1080302, 1080303, 1080304, 1080305, 1080306, 1080307, 1080308, 1080401, 1080402, 1080403, 1080404, 1080405, 1080406, 1080407, 1080408, 1080409, 1080410, 1080411, 1080412, 1080413, 1080414, 1080415, 1080416, 1080417, 1080418, 1080419, 1080420, 1080422, 1080426, 1090101, 1090102, 1090103, 1090104, 1090105, 1090106, 1090107, 1090109, 1090110, 1090111, 1090112, 1090113, 1090114, 1090115, 1100101, 1100102, 1100103, 1100104, 1100105, 1100106, 1100107, 1100108, 1100109, 1100110, 1100111, 1100112, 1100114, 1100115, 1100117, 1100118, 1100119, 1110133, 1110134, 1110135, 1110136, 1110137, 1110138, 1110139, 1110140, 1110141, 1110142, 1110143, 1110144, 1110145, 1110146, 1110147, 1110148, 1110149, 1110150, 1110151, 1110152, 1110153, 1110154, 1110156, 1110158, 1110159, 1110160, 1110201, 1110202, 1110203, 1110204, 1110205, 1110206, 1110207, 1110208, 1110209, 1110210, 1110211, 1110212, 1110213, 1110214, 1110215, 1110216, 1110217, 1110218, 1110219, 1110220, 1110221, 1110222, 1110223, 1110224, 1110225, 1110226, 1110227, 1110228, 1110229, 1110230, 1100116, 1070116, 1080421, 1080423, 1110301, 1110302, 1110303, 1110304, 1110305, 1110306, 1110307, 1110308, 1110309, 1110310, 1110311, 1110312, 1110313, 1110314, 1110315, 1110316, 1110317, 1110318, 1110319, 1110320, 1110321, 1110322, 1110323, 1110324, 1110325, 1110326)

for ($i = 0; $i -lt $filmIds.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $filmIds[$i]
}
